$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the "Date" value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-07-22T11:30:39+00:00"

# --- Mapping Table 0 sheet: rename "PersonneTierce" -> "Contact" in the
#     related-person rows (Source + Display columns, A & B) ---
$map0 = $wb.Worksheets.Item("Mapping Table 0")

$map0.Range("A3").Value = "Contact.IdContact"
$map0.Range("B3").Value = "Contact.IdContact"

$map0.Range("A4").Value = "Contact.adresse"
$map0.Range("B4").Value = "Contact.adresse"

$map0.Range("A5").Value = "Contact.telecommunication"
$map0.Range("B5").Value = "Contact.telecommunication"
